$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("2 jogos consecutivos)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.End - 1
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("TEMP_SPLIT", $bmRange)
$d.Bookmarks("TEMP_SPLIT").Delete()
# Now try inserting "=" right at pos
$insertPoint = $d.Range($pos, $pos)
$insertPoint.InsertBefore("=")
Write-Output "done"
